$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 546.3333
$ws.Range("I33").Value = 568.0909
$ws.Range("J33").Value = 486.5
$ws.Range("K33").Value = 568.0909
$ws.Range("L33").Value = 486.5
$ws.Range("M33").Value = -339.0909
$ws.Range("N33").Value = -944.5
$ws.Range("H62").Value = 13893114
$ws.Range("I62").Value = 15877416
$ws.Range("K62").Value = 15877416
$ws.Range("M62").Value = -15876792
$ws.Range("H65").Value = 13893114
$ws.Range("I65").Value = 15877416
$ws.Range("K65").Value = 79387080
$ws.Range("M65").Value = -79383960
$ws.Range("H96").Value = 1969.8
$ws.Range("I96").Value = 2384.1428
$ws.Range("K96").Value = 7152.428400000001
$ws.Range("M96").Value = -5779.428400000001
$ws.Range("H107").Value = 8252.546
$ws.Range("I107").Value = 7864.222
$ws.Range("K107").Value = 7864.222
$ws.Range("M107").Value = -5944.222
$ws.Range("H112").Value = 2047.45
$ws.Range("J112").Value = 2152.7222
$ws.Range("L112").Value = 6458.1666
$ws.Range("N112").Value = -8674.1666
$ws.Range("H125").Value = 2876.5908
$ws.Range("I125").Value = 3062.4443
$ws.Range("J125").Value = 2747.923
$ws.Range("K125").Value = 27561.9987
$ws.Range("L125").Value = 24731.307
$ws.Range("M125").Value = -25101.9987
$ws.Range("N125").Value = -29651.307
$ws.Range("H138").Value = 2791.5466
$ws.Range("I138").Value = 3099.4285
$ws.Range("J138").Value = 2759.853
$ws.Range("K138").Value = 9298.2855
$ws.Range("L138").Value = 8279.559000000001
$ws.Range("M138").Value = -4158.2855
$ws.Range("N138").Value = -18559.559

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11962.4
$ws.Range("I32").Value = 9189.648999999999
$ws.Range("J32").Value = 17606.928
$ws.Range("K32").Value = 9189.648999999999
$ws.Range("L32").Value = 17606.928
$ws.Range("M32").Value = -8902.648999999999
$ws.Range("N32").Value = -18180.928
$ws.Range("H61").Value = 52632920
$ws.Range("I61").Value = 71429390
$ws.Range("K61").Value = 71429390
$ws.Range("M61").Value = -71429178
$ws.Range("H74").Value = 2266.2144
$ws.Range("I74").Value = 1357
$ws.Range("K74").Value = 1357
$ws.Range("M74").Value = -483
$ws.Range("H77").Value = 2266.2144
$ws.Range("I77").Value = 1357
$ws.Range("K77").Value = 6785
$ws.Range("M77").Value = -2417
$ws.Range("H97").Value = 4600.8335
$ws.Range("I97").Value = 434.42856
$ws.Range("J97").Value = 33765.668
$ws.Range("K97").Value = 434.42856
$ws.Range("L97").Value = 33765.668
$ws.Range("M97").Value = 61.57144
$ws.Range("N97").Value = -34757.668
$ws.Range("H132").Value = 2768.2327
$ws.Range("I132").Value = 2357.8125
$ws.Range("J132").Value = 3962.182
$ws.Range("K132").Value = 7073.4375
$ws.Range("L132").Value = 11886.546
$ws.Range("M132").Value = -4543.4375
$ws.Range("N132").Value = -16946.546
$ws.Range("H136").Value = 52632920
$ws.Range("I136").Value = 71429390
$ws.Range("K136").Value = 214288170
$ws.Range("M136").Value = -214285620

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3643
$ws.Range("I20").Value = 3081.6
$ws.Range("J20").Value = 4344.75
$ws.Range("K20").Value = 3081.6
$ws.Range("L20").Value = 4344.75
$ws.Range("M20").Value = -2834.6
$ws.Range("N20").Value = -4838.75
$ws.Range("H86").Value = 34485508
$ws.Range("I86").Value = 41669188
$ws.Range("J86").Value = 3862.8
$ws.Range("K86").Value = 41669188
$ws.Range("L86").Value = 3862.8
$ws.Range("M86").Value = -41668065
$ws.Range("N86").Value = -6108.8
$ws.Range("H89").Value = 34485508
$ws.Range("I89").Value = 41669188
$ws.Range("J89").Value = 3862.8
$ws.Range("K89").Value = 208345940
$ws.Range("L89").Value = 19314
$ws.Range("M89").Value = -208340324
$ws.Range("N89").Value = -30546
$ws.Range("H94").Value = 13889852
$ws.Range("I94").Value = 17858082
$ws.Range("J94").Value = 1049.5
$ws.Range("K94").Value = 17858082
$ws.Range("L94").Value = 1049.5
$ws.Range("M94").Value = -17857631
$ws.Range("N94").Value = -1951.5
$ws.Range("H140").Value = 53304.75
$ws.Range("J140").Value = 53304.75
$ws.Range("L140").Value = 53304.75
$ws.Range("N140").Value = -63664.75
$ws.Range("H141").Value = 93260
$ws.Range("J141").Value = 93260
$ws.Range("L141").Value = 93260
$ws.Range("N141").Value = -103620

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 5169790
$ws.Range("I86").Value = 7432496
$ws.Range("K86").Value = 7432496
$ws.Range("M86").Value = -7431373
$ws.Range("H89").Value = 5169790
$ws.Range("I89").Value = 7432496
$ws.Range("K89").Value = 37162480
$ws.Range("M89").Value = -37156864
$ws.Range("H132").Value = 2032.0333
$ws.Range("I132").Value = 1679.2941
$ws.Range("J132").Value = 2493.3076
$ws.Range("K132").Value = 5037.8823
$ws.Range("L132").Value = 7479.9228
$ws.Range("M132").Value = -2507.8823
$ws.Range("N132").Value = -12539.9228
$ws.Range("H134").Value = 11112536
$ws.Range("I134").Value = 1450.742
$ws.Range("J134").Value = 35715656
$ws.Range("K134").Value = 4352.226
$ws.Range("L134").Value = 107146968
$ws.Range("M134").Value = -1817.226
$ws.Range("N134").Value = -107152038
$ws.Range("H141").Value = 254896
$ws.Range("J141").Value = 254896
$ws.Range("L141").Value = 254896
$ws.Range("N141").Value = -265256

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 16975494
$ws.Range("I131").Value = 111111560
$ws.Range("J131").Value = 31001.38
$ws.Range("K131").Value = 333334680
$ws.Range("L131").Value = 93004.14
$ws.Range("M131").Value = -333329640
$ws.Range("N131").Value = -103084.14

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 14027.5
$ws.Range("J57").Value = 23000
$ws.Range("L57").Value = 23000
$ws.Range("N57").Value = -24640
$ws.Range("H132").Value = 6308.8213
$ws.Range("I132").Value = 7706.722
$ws.Range("J132").Value = 3792.6
$ws.Range("K132").Value = 23120.166
$ws.Range("L132").Value = 11377.8
$ws.Range("M132").Value = -20590.166
$ws.Range("N132").Value = -16437.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2563
$ws.Range("I16").Value = 1414.3846
$ws.Range("J16").Value = 17495
$ws.Range("K16").Value = 1414.3846
$ws.Range("L16").Value = 17495
$ws.Range("M16").Value = -1244.3846
$ws.Range("N16").Value = -17835
$ws.Range("H100").Value = 1193.7
$ws.Range("I100").Value = 1079.75
$ws.Range("K100").Value = 1079.75
$ws.Range("M100").Value = -538.75
$ws.Range("H122").Value = 22729982
$ws.Range("I122").Value = 25002580
$ws.Range("K122").Value = 75007740
$ws.Range("M122").Value = -75005290
$ws.Range("H132").Value = 2559.4827
$ws.Range("I132").Value = 2128.1333
$ws.Range("K132").Value = 6384.3999
$ws.Range("M132").Value = -3854.3999
$ws.Range("H136").Value = 2153
$ws.Range("I136").Value = 2100.0625
$ws.Range("K136").Value = 6300.1875
$ws.Range("M136").Value = -3750.1875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 11366157
$ws.Range("I122").Value = 13891409
$ws.Range("J122").Value = 2526.25
$ws.Range("K122").Value = 41674227
$ws.Range("L122").Value = 7578.75
$ws.Range("M122").Value = -41671777
$ws.Range("N122").Value = -12478.75
$ws.Range("H132").Value = 4149.0557
$ws.Range("I132").Value = 4477.4644
$ws.Range("J132").Value = 2999.625
$ws.Range("K132").Value = 13432.3932
$ws.Range("L132").Value = 8998.875
$ws.Range("M132").Value = -10902.3932
$ws.Range("N132").Value = -14058.875
